# "started developing the csv-data import"
# The "Stand per ..." label (merged E3:F3) is shifted one column to the
# right so the freed-up column can later hold imported CSV data: the
# merge becomes F3:G3 and F3 now carries its own copy of the date label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F3 is currently the right half of the merged E3:F3 range, so it has to
# be unmerged before it can receive its own value.
$ws.Range("E3:F3").UnMerge()

# Give F3 the same text that already lives in E3.
$ws.Range("F3").Value = "Stand per 09.02.2018"

# Re-merge one column over: F3:G3 instead of the old E3:F3.
$ws.Range("F3:G3").Merge()
